$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "url"
$ws.Range("B1").Value = "status"
$ws.Range("A2").Value = "https://preview.allerganpro.com/co/es.html"
$ws.Range("B2").Value = 200
$ws.Range("A3").Value = "https://preview.allerganpro.com/libs/granite/csrf/token.json"
$ws.Range("B3").Value = 200
$ws.Range("A4").Value = "https://preview.allerganpro.com/bin/public/abbvie-commons/basic-login"
$ws.Range("B4").Value = 302
$ws.Range("A5").Value = "https://cag.abbvie.com:9999/bf/16a183f6-c871-4082-850b-a1f7a2ecd0b1?type=js3&sn=v_4_srv_-2D42_sn_FUIEST49VJHI6N83HHKMNSHTA28UU6KB&svrid=-42&flavor=cors&vi=MKHHAUQAKNCIRBPPMRQNSPFIQMMATKGB-0&modifiedSince=1665670355615&rf=https%3A%2F%2Fpreview.allerganpro.com%2Fbasic-login.html&bp=3&app=b90c0fbe356a6561&crc=1619779023&en=oao3vfhf&end=1"
$ws.Range("B5").Value = 200
$ws.Range("A6").Value = "https://cag.abbvie.com:9999/bf/16a183f6-c871-4082-850b-a1f7a2ecd0b1?type=js3&sn=v_4_srv_-2D42_sn_FUIEST49VJHI6N83HHKMNSHTA28UU6KB&svrid=-42&flavor=cors&vi=MKHHAUQAKNCIRBPPMRQNSPFIQMMATKGB-0&modifiedSince=1665670355615&rf=https%3A%2F%2Fpreview.allerganpro.com%2Fbasic-login.html&bp=3&app=b90c0fbe356a6561&crc=2003289253&en=oao3vfhf&end=1"
$ws.Range("B6").Value = 200
$ws.Range("A7").Value = "https://preview.allerganpro.com/etc.clientlibs/abbvie-pro/clientlibs/allergan-pro/publish-header.min.css"
$ws.Range("B7").Value = 200
$ws.Range("A8").Value = "https://preview.allerganpro.com/etc.clientlibs/abbvie-pro/clientlibs/allergan-pro/publish-header.min.js"
$ws.Range("B8").Value = 200
$ws.Range("A9").Value = "https://preview.allerganpro.com/etc.clientlibs/clientlibs/granite/jquery/granite/csrf.min.js"
$ws.Range("B9").Value = 200
$ws.Range("A10").Value = "https://preview.allerganpro.com/etc.clientlibs/abbvie-pro/components/content/headline-text/clientlibs.min.js"
$ws.Range("B10").Value = 200
$ws.Range("A11").Value = "https://preview.allerganpro.com/etc.clientlibs/abbvie-pro/components/content/button-link/clientlibs.min.css"
$ws.Range("B11").Value = 200
$ws.Range("A12").Value = "https://preview.allerganpro.com/etc.clientlibs/abbvie-pro/components/content/button-link/clientlibs.min.js"
$ws.Range("B12").Value = 200
$ws.Range("A13").Value = "https://cag.abbvie.com:9999/jstag/managed/ruxitagent_A2Vfqru_10249220905100923.js"
$ws.Range("B13").Value = 200
$ws.Range("A14").Value = "https://preview.allerganpro.com/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/fonts/roboto/roboto_bold.woff2"
$ws.Range("B14").Value = 200
$ws.Range("A15").Value = "https://preview.allerganpro.com/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/fonts/roboto/roboto_light.woff2"
$ws.Range("B15").Value = 200
$ws.Range("A16").Value = "https://preview.allerganpro.com/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/fonts/roboto/roboto_regular.woff2"
$ws.Range("B16").Value = 200
$ws.Range("A17").Value = "https://preview.allerganpro.com/co/es.html"
$ws.Range("B17").Value = 200
$ws.Range("A18").Value = "https://preview.allerganpro.com/etc.clientlibs/abbvie-pro/clientlibs/publish-footer.min.js"
$ws.Range("B18").Value = 200
$ws.Range("A19").Value = "https://preview.allerganpro.com/content/dam/abbvie-pro/co/abbvieprologo/AbbviePRO.png/_jcr_content/renditions/cq5dam.web.1280.1280.png"
$ws.Range("B19").Value = 200
$ws.Range("A20").Value = "https://preview.allerganpro.com/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/images/icons/Mobile%20nav.png"
$ws.Range("B20").Value = 200
$ws.Range("A21").Value = "https://consent.trustarc.com/v2/notice/hvz0wu"
$ws.Range("B21").Value = 200
$ws.Range("A22").Value = "https://preview.allerganpro.com/etc.clientlibs/abbvie-pro/components/content/image-extension/clientlibs.min.js"
$ws.Range("B22").Value = 200
$ws.Range("A23").Value = "https://preview.allerganpro.com/etc.clientlibs/abbvie-pro/components/content/meta-navigation/clientlibs.min.css"
$ws.Range("B23").Value = 200
$ws.Range("A24").Value = "https://preview.allerganpro.com/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/images/icons/mobile-menu-close.png"
$ws.Range("B24").Value = 200
$ws.Range("A25").Value = "https://preview.allerganpro.com/etc.clientlibs/abbvie-pro/components/content/meta-navigation/clientlibs.min.js"
$ws.Range("B25").Value = 200
$ws.Range("A26").Value = "https://consent.trustarc.com/v2/asset/ic-error.svg"
$ws.Range("B26").Value = 200
$ws.Range("A27").Value = "https://consent.trustarc.com/v2/asset/ic-close-white.svg"
$ws.Range("B27").Value = 200
$ws.Range("A28").Value = "https://consent.trustarc.com/v2/asset/ic-close.svg"
$ws.Range("B28").Value = 200
$ws.Range("A29").Value = "https://consent.trustarc.com/v2/asset/trustarc-logo-xs.svg"
$ws.Range("B29").Value = 200
$ws.Range("A30").Value = "https://consent.trustarc.com/v2/asset/latin.woff2"
$ws.Range("B30").Value = 200
$ws.Range("A31").Value = "https://preview.allerganpro.com/content/allergan-pro/co/es/jcr:content/header/header_area/image-extension/item_1.coreimg.png/1682019219951-AbbviePRO.png"
$ws.Range("B31").Value = 302
$ws.Range("A32").Value = "https://preview.allerganpro.com/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/fonts/hcpicon/hcpicon.ttf?q0neb3"
$ws.Range("B32").Value = 200
$ws.Range("A33").Value = "https://preview.allerganpro.com/libs/granite/csrf/token.json"
$ws.Range("B33").Value = 200
$ws.Range("A34").Value = "https://consent.trustarc.com/v2/consentcategories/getnonemptyindexes?cmId=hvz0wu&referer=&fullURL=https%3A%2F%2Fpreview.allerganpro.com%2Fco%2Fes.html&category="
$ws.Range("B34").Value = 200
$ws.Range("A35").Value = "https://consent-reporting.trustarc.com/api/user-action/log?action=impression&domain=hvz0wu&behavior=implied&country=bd&language=en&rand=0.5844605959308773&session=aa6cbbf4-d6fc-4622-a991-47ecda4fc6e3&userType=NEW"
$ws.Range("B35").Value = 202
$ws.Range("A36").Value = "https://preview.allerganpro.com/bin/public/abbvie-commons/hreflangs?resourcePath=/content/allergan-pro/co/es/jcr:content"
$ws.Range("B36").Value = 200
$ws.Range("A37").Value = "https://consent.trustarc.com/v2/asset/16:19:49.763hvz0wu_AbbVieID-logo.png"
$ws.Range("B37").Value = 200
$ws.Range("A38").Value = "https://preview.allerganpro.com/content/allergan-pro/co/es/jcr%3acontent/header/header_area/image-extension/item_1.coreimg.png/1707333003504.png"
$ws.Range("B38").Value = 200
$ws.Range("A39").Value = "https://consent-reporting.trustarc.com/api/user-action/bannermsg?action=views&domain=hvz0wu&behavior=implied&country=bd&language=en&rand=0.09455427329307686&session=aa6cbbf4-d6fc-4622-a991-47ecda4fc6e3&userType=NEW"
$ws.Range("B39").Value = 202
$ws.Range("A40").Value = "https://preview.allerganpro.com/content/dam/allergan-pro/colombia/home/Home_Articulos02.png/_jcr_content/renditions/cq5dam.web.1280.1280.png"
$ws.Range("B40").Value = 200
$ws.Range("A41").Value = "https://preview.allerganpro.com/content/dam/allergan-pro/colombia/home/Home_Articulos01.png/_jcr_content/renditions/cq5dam.web.1280.1280.png"
$ws.Range("B41").Value = 200
$ws.Range("A42").Value = "https://preview.allerganpro.com/content/dam/allergan-pro/colombia/home/New%20banner%20Home%20Allergan%20pro.png/_jcr_content/renditions/cq5dam.web.1280.1280.png"
$ws.Range("B42").Value = 200
$ws.Range("A43").Value = "https://preview.allerganpro.com/content/allergan-pro/co/es/jcr:content/body/column_control_copy/par1-100col/column_control_18337/par1-100col/column_control/par2-5050col/image_extension_copy/item_1.coreimg.png/1655224770198-Home_Articulos02.png"
$ws.Range("B43").Value = 200
$ws.Range("A44").Value = "https://preview.allerganpro.com/content/allergan-pro/co/es/jcr:content/body/column_control_copy/par1-100col/column_control_18337/par1-100col/column_control/par1-5050col/image_extension_copy/item_1.coreimg.png/1637251107888-Home_Articulos01.png"
$ws.Range("B44").Value = 200

$ws.Range("A45:B52").ClearContents()
